$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ08974452"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 5705.956255674148
$ws.Cells.Item(2, 3).Value = 0.006471042880512402
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 620.0786174113043
$ws.Cells.Item(3, 3).Value = 0.6937226549195461
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -767.0182632936534
$ws.Cells.Item(4, 3).Value = 0.5866451915609676
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -672.6010872226663
$ws.Cells.Item(5, 3).Value = 0.6334690357727504
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = -4281.182111053372
$ws.Cells.Item(6, 3).Value = [double]"1.762255662343111e-07"
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -3911.401961968803
$ws.Cells.Item(7, 3).Value = [double]"8.061567276866767e-07"
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 87.24374819528239
$ws.Cells.Item(8, 3).Value = 0.2346672639172685
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -732.0222202150586
$ws.Cells.Item(9, 3).Value = [double]"2.849804689972911e-05"
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -12.6508528578513
$ws.Cells.Item(10, 3).Value = 0.08555537900937248
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = 188.0031570827565
$ws.Cells.Item(11, 3).Value = 0.05392927108304432
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 946.196651871971
$ws.Cells.Item(12, 3).Value = [double]"3.73016558035957e-42"
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = 0.0009153130696142736
$ws.Cells.Item(13, 3).Value = 0.9919530804930832
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = [double]"7.194548136279819e-05"
$ws.Cells.Item(14, 3).Value = 0.4689258349414102
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = 2.110680647264893
$ws.Cells.Item(15, 3).Value = 0.8137750584710559
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 13.71068305218201
$ws.Cells.Item(16, 3).Value = 0.1153474773521193
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -1685.891628989182
$ws.Cells.Item(17, 3).Value = 0.1022206548160452
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = -243.7067556322393
$ws.Cells.Item(18, 3).Value = 0.7912147820858623
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 2875.596182427542
$ws.Cells.Item(19, 3).Value = 0.3579970034156897

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ09194252"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 2971.204119603321
$ws.Cells.Item(2, 3).Value = 0.1418826045353619
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 517.4803710481592
$ws.Cells.Item(3, 3).Value = 0.7425992082251018
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -558.7966263894505
$ws.Cells.Item(4, 3).Value = 0.6849144672838061
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -593.7827764090453
$ws.Cells.Item(5, 3).Value = 0.6665600329848409
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = -1011.936546107465
$ws.Cells.Item(6, 3).Value = 0.1690109546229872
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -772.9623071142133
$ws.Cells.Item(7, 3).Value = 0.2742081020771847
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 33.72138699383011
$ws.Cells.Item(8, 3).Value = 0.6361564783083677
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -905.4841842993859
$ws.Cells.Item(9, 3).Value = [double]"1.557385622905376e-07"
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -15.94186116341721
$ws.Cells.Item(10, 3).Value = 0.02633095817427787
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = 324.404063098086
$ws.Cells.Item(11, 3).Value = 0.0008234352268144832
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 874.3301193433604
$ws.Cells.Item(12, 3).Value = [double]"7.285006534023662e-39"
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = -0.01297374000572531
$ws.Cells.Item(13, 3).Value = 0.8831295529731762
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = [double]"8.944369636368953e-05"
$ws.Cells.Item(14, 3).Value = 0.3555226543508384
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -1.759151497702486
$ws.Cells.Item(15, 3).Value = 0.8387912215525581
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 11.7612167007617
$ws.Cells.Item(16, 3).Value = 0.1656476295383076
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -1747.806791376361
$ws.Cells.Item(17, 3).Value = 0.07347815025993081
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = 147.9015717407184
$ws.Cells.Item(18, 3).Value = 0.8731799701069407
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 1177.678110712885
$ws.Cells.Item(19, 3).Value = 0.6976076766614725

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ09422416"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 2238.288810658404
$ws.Cells.Item(2, 3).Value = 0.28272577374078
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 723.8200851535037
$ws.Cells.Item(3, 3).Value = 0.6513555087146595
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -743.7251363916969
$ws.Cells.Item(4, 3).Value = 0.599820717784738
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -524.4674365395103
$ws.Cells.Item(5, 3).Value = 0.7116417033405259
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = -640.391611579772
$ws.Cells.Item(6, 3).Value = 0.374615923743343
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -339.7504817510485
$ws.Cells.Item(7, 3).Value = 0.622570653996098
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 95.54704016395098
$ws.Cells.Item(8, 3).Value = 0.1929561183996466
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -658.4037230193858
$ws.Cells.Item(9, 3).Value = 0.0001994666915021765
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -14.7605778383804
$ws.Cells.Item(10, 3).Value = 0.04740952803896632
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = 159.8837852615528
$ws.Cells.Item(11, 3).Value = 0.0995078996944838
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 812.6144485869727
$ws.Cells.Item(12, 3).Value = [double]"1.07857658959971e-31"
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = 0.05916103884301432
$ws.Cells.Item(13, 3).Value = 0.5178038430516134
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = [double]"3.601408938691423e-05"
$ws.Cells.Item(14, 3).Value = 0.7172313236722963
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -6.797258992145732
$ws.Cells.Item(15, 3).Value = 0.4507043489986007
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 21.15152215323148
$ws.Cells.Item(16, 3).Value = 0.01669387848025958
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -712.884030232795
$ws.Cells.Item(17, 3).Value = 0.4771209240832703
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = -125.5869195313855
$ws.Cells.Item(18, 3).Value = 0.8927039344573175
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = -648.6898319159345
$ws.Cells.Item(19, 3).Value = 0.8363497050999875

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ09627437"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = -962.9180500736556
$ws.Cells.Item(2, 3).Value = 0.7019883216041874
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 3210.295122532678
$ws.Cells.Item(3, 3).Value = 0.1467087714165947
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = 1657.122604676289
$ws.Cells.Item(4, 3).Value = 0.4298110567506942
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 1858.945641878318
$ws.Cells.Item(5, 3).Value = 0.3758478500070834
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = -663.8197434717466
$ws.Cells.Item(6, 3).Value = 0.3648589615622777
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -353.4330452659489
$ws.Cells.Item(7, 3).Value = 0.6153880640723609
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 62.6813161835085
$ws.Cells.Item(8, 3).Value = 0.3855501118050445
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -459.8794917949842
$ws.Cells.Item(9, 3).Value = 0.008038400459471637
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -22.10363905486413
$ws.Cells.Item(10, 3).Value = 0.002806469372078432
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = 188.7666076678642
$ws.Cells.Item(11, 3).Value = 0.04992304659365807
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 915.4740381720353
$ws.Cells.Item(12, 3).Value = [double]"8.147862662226763e-39"
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = 0.04546613657323023
$ws.Cells.Item(13, 3).Value = 0.6078822256792855
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = [double]"2.257985185900844e-05"
$ws.Cells.Item(14, 3).Value = 0.8175108511320086
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -5.389294095188369
$ws.Cells.Item(15, 3).Value = 0.5441652364452059
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 22.54493417945308
$ws.Cells.Item(16, 3).Value = 0.009627624198253304
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -298.9257700531919
$ws.Cells.Item(17, 3).Value = 0.761263702990181
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = 577.6994594343406
$ws.Cells.Item(18, 3).Value = 0.5261703833946338
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 3812.482078886899
$ws.Cells.Item(19, 3).Value = 0.2151853457701389

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ09823589"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 1087.04462038417
$ws.Cells.Item(2, 3).Value = 0.6075349063289598
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 1240.343008236481
$ws.Cells.Item(3, 3).Value = 0.4409883429282825
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -390.2164637107901
$ws.Cells.Item(4, 3).Value = 0.7840844140619847
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -117.9411639136122
$ws.Cells.Item(5, 3).Value = 0.9340572861678477
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = -1710.141283455174
$ws.Cells.Item(6, 3).Value = 0.01507679082960434
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -1186.419211899391
$ws.Cells.Item(7, 3).Value = 0.07681372466204488
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 24.21281098292347
$ws.Cells.Item(8, 3).Value = 0.7428524876779038
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -766.2741028002386
$ws.Cells.Item(9, 3).Value = [double]"1.463201795771607e-05"
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -15.81168820669041
$ws.Cells.Item(10, 3).Value = 0.03164582706258822
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = 205.9176238682899
$ws.Cells.Item(11, 3).Value = 0.03650954272789097
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 889.3280014265523
$ws.Cells.Item(12, 3).Value = [double]"1.412351434157667e-35"
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = 0.04638940149675023
$ws.Cells.Item(13, 3).Value = 0.6110093595304464
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = [double]"3.075360182235486e-05"
$ws.Cells.Item(14, 3).Value = 0.7593529643464272
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -1.169757551273947
$ws.Cells.Item(15, 3).Value = 0.8967922535944275
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 31.91658879973383
$ws.Cells.Item(16, 3).Value = 0.0004288915680245612
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -947.486088634084
$ws.Cells.Item(17, 3).Value = 0.3411001684423538
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = 255.1315202180974
$ws.Cells.Item(18, 3).Value = 0.7854224589076002
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 3132.0648207863
$ws.Cells.Item(19, 3).Value = 0.3183983784983765

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ10032611"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 3603.640888172667
$ws.Cells.Item(2, 3).Value = 0.09152378930328899
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 162.0598886505243
$ws.Cells.Item(3, 3).Value = 0.9238422686113559
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -1988.572170424025
$ws.Cells.Item(4, 3).Value = 0.188531801926698
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -1805.962222272133
$ws.Cells.Item(5, 3).Value = 0.2321874111127512
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = -1565.615715173758
$ws.Cells.Item(6, 3).Value = 0.0300714183107951
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -1378.744011527601
$ws.Cells.Item(7, 3).Value = 0.04601188443214888
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 36.14227613250988
$ws.Cells.Item(8, 3).Value = 0.6114648569770784
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -880.7271942604651
$ws.Cells.Item(9, 3).Value = [double]"3.91454820180948e-07"
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -16.99940864072073
$ws.Cells.Item(10, 3).Value = 0.02041289737563938
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = 84.09936802209955
$ws.Cells.Item(11, 3).Value = 0.3849288136480108
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 849.21688459302
$ws.Cells.Item(12, 3).Value = [double]"1.160983815632575e-34"
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = 0.01132754755919805
$ws.Cells.Item(13, 3).Value = 0.8991255448633657
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = [double]"2.728317658862452e-05"
$ws.Cells.Item(14, 3).Value = 0.7845943387623728
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = 1.734287520279237
$ws.Cells.Item(15, 3).Value = 0.8449377926344546
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 30.37256984850789
$ws.Cells.Item(16, 3).Value = 0.0004682910215742239
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -1548.198458110263
$ws.Cells.Item(17, 3).Value = 0.1415964426412777
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = 295.4990738120262
$ws.Cells.Item(18, 3).Value = 0.7453950258858705
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 1892.820472745905
$ws.Cells.Item(19, 3).Value = 0.5379455513146865

# --- Sheet 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ10238631"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 1667.492402649698
$ws.Cells.Item(2, 3).Value = 0.4244270190018434
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 1048.625443435974
$ws.Cells.Item(3, 3).Value = 0.5083084125134059
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -410.6755150680565
$ws.Cells.Item(4, 3).Value = 0.7724505265833971
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -314.2304509603099
$ws.Cells.Item(5, 3).Value = 0.8249145123459607
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = -1103.279628957597
$ws.Cells.Item(6, 3).Value = 0.1395667026553119
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -845.4095598750021
$ws.Cells.Item(7, 3).Value = 0.2363536839079075
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 57.04725348956116
$ws.Cells.Item(8, 3).Value = 0.4367405147366987
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -777.395069874953
$ws.Cells.Item(9, 3).Value = [double]"1.136459403068022e-05"
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -20.69203027041258
$ws.Cells.Item(10, 3).Value = 0.005643880988182599
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = 185.8417595899911
$ws.Cells.Item(11, 3).Value = 0.05904098459086143
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 883.1322766060573
$ws.Cells.Item(12, 3).Value = [double]"1.175801260341223e-36"
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = -0.03724023592930122
$ws.Cells.Item(13, 3).Value = 0.6860396641638149
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = 0.0001003433229585861
$ws.Cells.Item(14, 3).Value = 0.320680610799756
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -2.918621700339628
$ws.Cells.Item(15, 3).Value = 0.7469564244626271
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 26.07729400291375
$ws.Cells.Item(16, 3).Value = 0.003220390266651215
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -1041.501766631232
$ws.Cells.Item(17, 3).Value = 0.3147255502463364
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = 847.7361116563076
$ws.Cells.Item(18, 3).Value = 0.3820145824525514
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 3224.816664526028
$ws.Cells.Item(19, 3).Value = 0.3077266935125695

# --- Sheet 8 ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ10440655"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 1512.510859298831
$ws.Cells.Item(2, 3).Value = 0.5093109343207751
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 722.6962128795117
$ws.Cells.Item(3, 3).Value = 0.6979832953594959
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -258.821226950251
$ws.Cells.Item(4, 3).Value = 0.8770714243414196
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = 102.945731185534
$ws.Cells.Item(5, 3).Value = 0.951027901364585
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = -1874.922295491788
$ws.Cells.Item(6, 3).Value = 0.01072768914057577
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -1509.825713571014
$ws.Cells.Item(7, 3).Value = 0.03231782058222307
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 12.3897395949798
$ws.Cells.Item(8, 3).Value = 0.8647141788242829
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -760.3486753931454
$ws.Cells.Item(9, 3).Value = [double]"1.388587591655887e-05"
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -16.00807570105419
$ws.Cells.Item(10, 3).Value = 0.02879664440982044
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = 190.4267672319682
$ws.Cells.Item(11, 3).Value = 0.05028967122639668
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 905.782962973866
$ws.Cells.Item(12, 3).Value = [double]"6.163740775106954e-39"
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = 0.06774294156149711
$ws.Cells.Item(13, 3).Value = 0.452467137263778
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = [double]"4.213556529219741e-05"
$ws.Cells.Item(14, 3).Value = 0.669638879800174
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = 7.399340972372105
$ws.Cells.Item(15, 3).Value = 0.4003229524017605
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 28.46206248319082
$ws.Cells.Item(16, 3).Value = 0.0009365603878591209
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -1324.848913470365
$ws.Cells.Item(17, 3).Value = 0.1863802374765029
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = -703.7519419369299
$ws.Cells.Item(18, 3).Value = 0.4478988899331033
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = -1414.937066547522
$ws.Cells.Item(19, 3).Value = 0.6478840456711089

# --- Sheet 9 ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ10644376"
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 3424.358438459255
$ws.Cells.Item(2, 3).Value = 0.08893488980092878
$ws.Cells.Item(3, 1).Value = "Education[T.Primary/None]"
$ws.Cells.Item(3, 2).Value = 1180.011708696253
$ws.Cells.Item(3, 3).Value = 0.4544518496756211
$ws.Cells.Item(4, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(4, 2).Value = -868.565306909391
$ws.Cells.Item(4, 3).Value = 0.5290341366768405
$ws.Cells.Item(5, 1).Value = "Education[T.University]"
$ws.Cells.Item(5, 2).Value = -591.0814749409265
$ws.Cells.Item(5, 3).Value = 0.6686750278618472
$ws.Cells.Item(6, 1).Value = "Season[T.Spring]"
$ws.Cells.Item(6, 2).Value = -1379.65256250362
$ws.Cells.Item(6, 3).Value = 0.0430854944409556
$ws.Cells.Item(7, 1).Value = "Season[T.Winter]"
$ws.Cells.Item(7, 2).Value = -1034.686053383921
$ws.Cells.Item(7, 3).Value = 0.1117905666380415
$ws.Cells.Item(8, 1).Value = "HHSize"
$ws.Cells.Item(8, 2).Value = 36.5477667193006
$ws.Cells.Item(8, 3).Value = 0.6145969327192206
$ws.Cells.Item(9, 1).Value = "Sex"
$ws.Cells.Item(9, 2).Value = -819.0337479000592
$ws.Cells.Item(9, 3).Value = [double]"1.93662285009422e-06"
$ws.Cells.Item(10, 1).Value = "Age"
$ws.Cells.Item(10, 2).Value = -10.26868764317417
$ws.Cells.Item(10, 3).Value = 0.1590832993101499
$ws.Cells.Item(11, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(11, 2).Value = 71.11448177064489
$ws.Cells.Item(11, 3).Value = 0.4640277001596954
$ws.Cells.Item(12, 1).Value = "DistCenter_res"
$ws.Cells.Item(12, 2).Value = 890.5417767711859
$ws.Cells.Item(12, 3).Value = [double]"8.721994427771798e-39"
$ws.Cells.Item(13, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(13, 2).Value = -0.0008579651184658782
$ws.Cells.Item(13, 3).Value = 0.9923141302341066
$ws.Cells.Item(14, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(14, 2).Value = [double]"9.949975564107549e-05"
$ws.Cells.Item(14, 3).Value = 0.310668988881326
$ws.Cells.Item(15, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(15, 2).Value = -6.76882509481731
$ws.Cells.Item(15, 3).Value = 0.4385776239953125
$ws.Cells.Item(16, 1).Value = "street_length_res"
$ws.Cells.Item(16, 2).Value = 17.27641167944755
$ws.Cells.Item(16, 3).Value = 0.04536667407407577
$ws.Cells.Item(17, 1).Value = "LU_Comm_res"
$ws.Cells.Item(17, 2).Value = -2822.201592565309
$ws.Cells.Item(17, 3).Value = 0.006327824691041529
$ws.Cells.Item(18, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(18, 2).Value = 306.3631386053827
$ws.Cells.Item(18, 3).Value = 0.7360862503190033
$ws.Cells.Item(19, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(19, 2).Value = 4139.479973686363
$ws.Cells.Item(19, 3).Value = 0.1764851383272801

